$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to Text so numeric-looking strings (e.g. "1.010",
    # "26.619.25") are stored verbatim instead of being parsed as numbers,
    # then strip the temporary number-format override so the cell style
    # index is left exactly as it was before the write.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "26.631.65"
$ws.Range("E2").Value = "  +1.14%  "
Set-TextValue "D3" "1.825.93"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.58%  "
Set-TextValue "D5" "1.008"
$ws.Range("E5").Value = "  +0.45%  "
Set-TextValue "D6" "308.56"
$ws.Range("E6").Value = "  +0.59%  "
Set-TextValue "D7" "0.4678"
$ws.Range("E7").Value = "  +3.82%  "
Set-TextValue "D8" "0.3606"
$ws.Range("E8").Value = "  +0.37%  "
Set-TextValue "D9" "0.07132"
$ws.Range("E9").Value = "  +0.69%  "
Set-TextValue "D10" "0.9023"
$ws.Range("E10").Value = "  +2.09%  "
Set-TextValue "D11" "0.07752"
$ws.Range("E11").Value = "  +0.05%  "
Set-TextValue "D12" "19.42"
$ws.Range("E12").Value = "  -0.08%  "
Set-TextValue "D13" "1.825.46"
$ws.Range("E13").Value = "  +1.99%  "
Set-TextValue "D14" "5.272"
$ws.Range("E14").Value = "  -0.06%  "
Set-TextValue "D15" "6.355"
$ws.Range("E15").Value = "  +0.52%  "
Set-TextValue "D16" "87.53"
$ws.Range("E16").Value = "  +3.18%  "
Set-TextValue "D17" "1.010"
$ws.Range("E17").Value = "  +0.46%  "
Set-TextValue "D18" "0.000008549"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.35%  "
Set-TextValue "D20" "26.665.95"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  -0.34%  "
Set-TextValue "D22" "5.023"
$ws.Range("E22").Value = "  +1.05%  "
Set-TextValue "D23" "10.55"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -3.51%  "
Set-TextValue "D25" "152.98"
$ws.Range("E25").Value = "  +1.03%  "
Set-TextValue "D26" "17.92"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -2.09%  "
Set-TextValue "D28" "113.89"
$ws.Range("E28").Value = "  +1.78%  "
Set-TextValue "D29" "4.872"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("E32").Value = "  +3.16%  "
Set-TextValue "D33" "1.166"
$ws.Range("E33").Value = "  +5.46%  "
Set-TextValue "D34" "0.7359"
$ws.Range("E34").Value = "  +1.77%  "
Set-TextValue "D35" "4.441"
$ws.Range("E35").Value = "  -0.04%  "
Set-TextValue "D36" "1.080"
Set-TextValue "D37" "0.01930"
$ws.Range("E37").Value = "  +0.12%  "
Set-TextValue "D38" "0.05160"
$ws.Range("E38").Value = "  +1.45%  "
Set-TextValue "D39" "2.900"
$ws.Range("E39").Value = "  +1.64%  "
Set-TextValue "D40" "6.874"
$ws.Range("E40").Value = "  +0.36%  "
Set-TextValue "D41" "0.5053"
$ws.Range("E41").Value = "  -0.31%  "
Set-TextValue "D42" "0.1496"
$ws.Range("E42").Value = "  -1.24%  "
Set-TextValue "D43" "8.035"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  +0.73%  "
Set-TextValue "D46" "10.02"
$ws.Range("E46").Value = "  +2.02%  "
Set-TextValue "D47" "97.98"
$ws.Range("E47").Value = "  -3.12%  "
Set-TextValue "D48" "1.573"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +1.58%  "
Set-TextValue "D50" "64.01"
$ws.Range("E50").Value = "  -0.15%  "
Set-TextValue "D51" "35.78"
$ws.Range("E51").Value = "  -0.38%  "
